$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "7:15PM 12-30-2017"
$ws.Range("C15").Value = 74

$ws.Range("A16").Value = "10:00AM 12-31-2017"
$ws.Range("B16").Value = "12:52PM 12-31-2017"
$ws.Range("C16").Value = 172

$ws.Range("C17").Formula = "=SUM(C2:C16)/60"

$ws.Range("C18").Select()
